$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time value in B5 (Jacob's time) from 1.5 to 3.5
$ws.Range("B5").Value = 3.5

# Move/update the active selection to E9
$ws.Range("E9").Select()
